$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.650.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.327.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.385"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.83%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.960"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.327.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.194"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.390.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.942.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000243"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.324.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "492.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.452"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000184"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.497.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.53%  "

$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.139"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.60"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.96%  "

$ws.Range("E34").Value = "  +0.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.171"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.527"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "556.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.870"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0410"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.17%  "
